$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 4 (for the Dagger_Hit particle entity) - this pushes the
# existing Axe_1_Epic_Skill row (old row 4) down to row 5.
$ws.Rows.Item(4).Insert()

# Insert a new row 6 (for the Sword_1 Epic Skill entity), placed after the
# Axe_1_Epic_Skill row which is now row 5.
$ws.Rows.Item(6).Insert()

# Fill row 6 first so its string lands in the shared string table before the
# row 4 string (matches authoring order of the target workbook).
$ws.Range("A6").Value = 10112014
$ws.Range("B6").Value = "Items/Prefabs/Weapons/Skill/Sword_1_Epic_Skill"

# New row for the Dagger_Hit particle entity.
$ws.Range("A4").Value = 71000003
$ws.Range("B4").Value = "Particle/Prefabs/Dagger_Hit"

# Row 6 inherited the highlighted style from row 5 on insert; the new row
# should be unstyled like the other freshly-added rows.
$ws.Range("A6").ClearFormats()

# Restore the selected cell shown in the edited workbook.
$ws.Range("D17").Select()
